$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.492.54"
$ws.Range("E2").Value = "  +1.68%  "

$ws.Range("D3").Value = "2.326.57"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Formula = "'541.37"
$ws.Range("E5").Value = "  +5.30%  "

$ws.Range("D6").Formula = "'134.71"
$ws.Range("E6").Value = "  +1.63%  "

$ws.Range("D7").Formula = "'0.994"
$ws.Range("E7").Value = "  -0.49%  "

$ws.Range("D8").Formula = "'0.537"
$ws.Range("E8").Value = "  +0.64%  "

$ws.Range("D9").Value = "2.356.02"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("D10").Formula = "'0.103"
$ws.Range("E10").Value = "  +1.72%  "

$ws.Range("E11").Value = "  +0.89%  "

$ws.Range("E12").Value = "  +1.81%  "

$ws.Range("E13").Value = "  +4.38%  "

$ws.Range("D14").Value = "2.776.25"
$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").Formula = "'23.50"
$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").Value = "57.572.06"
$ws.Range("E16").Value = "  +1.84%  "

$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").Value = "2.332.64"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("D20").Formula = "'334.87"
$ws.Range("E20").Value = "  +2.85%  "

$ws.Range("E21").Value = "  +1.49%  "

$ws.Range("D22").Formula = "'6.73"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").Formula = "'1.00"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").Formula = "'61.60"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("E25").Value = "  +3.87%  "

$ws.Range("D26").Formula = "'0.997"
$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").Formula = "'8.40"
$ws.Range("E27").Value = "  -4.15%  "

$ws.Range("D28").Formula = "'1.42"
$ws.Range("E28").Value = "  +8.53%  "

$ws.Range("E29").Value = "  +4.40%  "

$ws.Range("D30").Formula = "'170.27"
$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("D31").Value = "0.0₃0737"
$ws.Range("E31").Value = "  +2.04%  "

$ws.Range("D32").Formula = "'6.17"
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("D33").Formula = "'18.56"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("E34").Value = "  +15.01%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").Formula = "'0.991"
$ws.Range("E36").Value = "  -0.58%  "

$ws.Range("E37").Value = "  -0.42%  "

$ws.Range("E38").Value = "  +5.20%  "

$ws.Range("E39").Value = "  +2.77%  "

$ws.Range("D40").Formula = "'39.35"
$ws.Range("E40").Value = "  +2.28%  "

$ws.Range("D41").Formula = "'150.94"
$ws.Range("E41").Value = "  -1.83%  "

$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  +1.15%  "

$ws.Range("D44").Formula = "'285.53"
$ws.Range("E44").Value = "  +1.93%  "

$ws.Range("D45").Formula = "'19.33"
$ws.Range("E45").Value = "  +5.81%  "

$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("E47").Value = "  +1.87%  "

$ws.Range("E48").Value = "  +0.46%  "

$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("D50").Formula = "'17.58"
$ws.Range("E50").Value = "  +1.86%  "

$ws.Range("D51").Formula = "'0.383"
$ws.Range("E51").Value = "  +0.17%  "
